# Apply the "Anonymize fedcore" edits described in the commit:
#  - rename every "fedcore" header label to "approach"
#  - add a thin top/bottom(/right) border under the header row for the
#    empty cells that sit beside the merged "original/fedcore/change" title
#  - normalize the "-0" delta values to plain "0"
#  - drop the stray empty cell at G5 on the computational_comparison sheet

$wb = $excel.ActiveWorkbook

# --- Sheet 1: quality_comparison -------------------------------------
$ws1 = $wb.Worksheets.Item("quality_comparison")

$ws1.Range("C2").Value = "approach"

# Build the "top+bottom" bordered style on C1 from a clean slate...
$c1 = $ws1.Range("C1")
$c1.ClearFormats()
$c1.Borders.Item(8).LineStyle = 1   # xlEdgeTop
$c1.Borders.Item(9).LineStyle = 1   # xlEdgeBottom

# ...then clone it onto D1 (format-only paste, so D1 inherits the exact
# same underlying style record instead of forking its own) and extend it
# with the extra right edge to get the "top+bottom+right" style.
$d1 = $ws1.Range("D1")
$c1.Copy()
$d1.PasteSpecial(-4122)             # xlPasteFormats
$d1.Borders.Item(10).LineStyle = 1  # xlEdgeRight

$ws1.Range("D4").Value = 0
$ws1.Range("D5").Value = 0
$ws1.Range("D12").Value = 0

# --- Sheet 2: computational_comparison --------------------------------
$ws2 = $wb.Worksheets.Item("computational_comparison")

$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# Reuse the two styles already built on sheet 1 (format-only paste keeps
# the workbook-level style table from growing / forking again).
$c1b = $ws2.Range("C1")
$c1.Copy()
$c1b.PasteSpecial(-4122)

$d1b = $ws2.Range("D1")
$d1.Copy()
$d1b.PasteSpecial(-4122)

$f1 = $ws2.Range("F1")
$c1.Copy()
$f1.PasteSpecial(-4122)

$g1 = $ws2.Range("G1")
$d1.Copy()
$g1.PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Drop the stray empty cell left over at G5
$ws2.Range("G5").ClearContents()
